# Fix misspelled shared string: "THree-Point Attempts" -> "Three-Point Attempts"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q35").Value = "Three-Point Attempts"

# ---------------------------------------------------------------------------
# Styling pass. Three visual roles are introduced:
#   - "title"  : section title cells (col A only)        -> bold, 14pt
#   - "header" : column-header rows                       -> bold, 14pt, red fill, thin border
#   - "data"   : the data rows beneath each header         -> regular, blue fill, thin border
# plus a secondary, unrelated 1-column mini table (Q31:Q43) that reuses the
# same red/bold "header" look for its own header cell, and a dedicated green
# highlight for its body cells.
# ---------------------------------------------------------------------------

function Style-Title($cell) {
    $cell.Font.Bold = $true
    $cell.Font.Size = 14
}

function Style-Header($rng) {
    $rng.Font.Bold = $true
    $rng.Font.Size = 14
    $rng.Interior.Color = 8421631        # RGB(255,128,128)
    $rng.Interior.PatternColor = 8421631
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
}

function Style-HeaderNoBorder($rng) {
    $rng.Font.Bold = $true
    $rng.Font.Size = 14
    $rng.Interior.Color = 8421631        # RGB(255,128,128)
    $rng.Interior.PatternColor = 8421631
}

function Style-Data($rng) {
    $rng.Interior.Color = 16767411       # RGB(179,217,255) / hex b3d9ff
    $rng.Interior.PatternColor = 16767411
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
}

function Style-DataNoBorder($rng) {
    $rng.Interior.Color = 9240346        # RGB(26,255,140) / hex 1aff8c
    $rng.Interior.PatternColor = 9240346
}

# --- Title cells (column A, standalone section titles) ---------------------
$titleRows = @(1,4,8,12,16,20,24)
foreach ($r in $titleRows) {
    Style-Title ($ws.Cells.Item($r, 1))
}

# --- The five (header,data) row pairs sharing the same layout --------------
Style-Header ($ws.Range("A5:X5"))
Style-Data   ($ws.Range("A6:X6"))

Style-Header ($ws.Range("A9:X9"))
Style-Data   ($ws.Range("A10:X10"))

Style-Header ($ws.Range("A13:X13"))
Style-Data   ($ws.Range("A14:X14"))

Style-Header ($ws.Range("A17:X17"))
Style-Data   ($ws.Range("A18:X18"))

Style-Header ($ws.Range("A21:U21"))
Style-Data   ($ws.Range("A22:U22"))

Style-Header ($ws.Range("A25:T25"))
Style-Data   ($ws.Range("A26:T26"))

# --- Play type table: header row, then a labeled column + data block -------
Style-Header ($ws.Range("A30:O30"))
Style-Header ($ws.Range("A31:A41"))
Style-Data   ($ws.Range("B31:O41"))

# --- Secondary 1-column mini table (stat-name legend) in column Q ----------
Style-HeaderNoBorder ($ws.Range("Q31"))
Style-DataNoBorder   ($ws.Range("Q32:Q43"))

Write-Host "styling complete"
